$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 106
$ws.Range("F3").Value = 179
$ws.Range("F4").Value = 416
$ws.Range("F5").Value = 195
$ws.Range("F6").Value = 132
$ws.Range("F7").Value = 1141
$ws.Range("F8").Value = 383
$ws.Range("F9").Value = 194
$ws.Range("F10").Value = 51
$ws.Range("F12").Value = 376
$ws.Range("F13").Value = 396
$ws.Range("F14").Value = 786
$ws.Range("F15").Value = 176
$ws.Range("F16").Value = 721
$ws.Range("F17").Value = 284
$ws.Range("F19").Value = 1009
$ws.Range("F20").Value = 457
$ws.Range("F21").Value = 264
$ws.Range("F23").Value = 380
$ws.Range("F25").Value = 41

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 284
$ws.Range("F9").Value = 9

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 345

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 345
$ws.Range("F4").Value = 106
$ws.Range("F5").Value = 179
$ws.Range("F6").Value = 416
$ws.Range("F7").Value = 195
$ws.Range("F8").Value = 132
$ws.Range("F9").Value = 1141
$ws.Range("F10").Value = 383
$ws.Range("F11").Value = 194
$ws.Range("F13").Value = 51
$ws.Range("F17").Value = 376
$ws.Range("F19").Value = 284
$ws.Range("F20").Value = 396
$ws.Range("F21").Value = 786
$ws.Range("F22").Value = 176
$ws.Range("F23").Value = 721
$ws.Range("F24").Value = 284
$ws.Range("F26").Value = 1009
$ws.Range("F27").Value = 457
$ws.Range("F29").Value = 9
$ws.Range("F30").Value = 264
$ws.Range("F32").Value = 380
$ws.Range("F36").Value = 41
